$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.244.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "'2.308.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'301.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "'98.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").Value = "'0.518"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.31%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "'17.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("D14").Value = "'6.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "'2.666.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "'2.308.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "'0.792"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "'43.100.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").Value = "'13.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.65%  "
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'68.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").Value = "'238.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "'2.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "'25.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("D29").Value = "'166.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'9.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("E31").Value = "  -6.45%  "
$ws.Range("D32").Value = "'33.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("D33").Value = "'5.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'18.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.15%  "
$ws.Range("D36").Value = "'4.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.33%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").Value = "'0.103"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("D43").Value = "'2.013.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("E45").Value = "  -6.96%  "
$ws.Range("D46").Value = "'10.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("D47").Value = "'17.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").Value = "'54.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").Value = "'2.538.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("E51").Value = "  +0.80%  "
